$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing event name re-typed in uppercase ---
$ws.Range("B133").Value = "HOPSEN BY AREA51"

# --- Row 134: new event (ATHÉNA @ Stahlwerk, Düsseldorf) ---
$ws.Range("A134").Value = 45752
$ws.Range("B134").Value = "ATHÉNA"
$ws.Range("B134").NumberFormat = "@"
$ws.Range("C134").Value = "Stahlwerk"
$ws.Range("C134").NumberFormat = "@"
$ws.Range("D134").Value = "Düsseldorf"
$ws.Range("D134").NumberFormat = "@"

$ws.Range("E134").Value = "https://www.instagram.com/reel/DFnmlijMBzP/?igsh=NnV2NGpsNWF6NDls"
$ws.Range("E134").NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Range("E134"), "https://www.instagram.com/reel/DFnmlijMBzP/?igsh=NnV2NGpsNWF6NDls", "", "", "https://www.instagram.com/reel/DFnmlijMBzP/?igsh=NnV2NGpsNWF6NDls") | Out-Null

$e134a = $ws.Range("E134").Characters(1, 64)
$e134a.Font.Underline = $true
$e134a.Font.Color = 16711680
$e134b = $ws.Range("E134").Characters(65, 1)
$e134b.Font.Underline = $true
$e134b.Font.Color = 16711680

$ws.Range("E134").Font.Name = "Calibri"
$ws.Range("E134").Font.Size = 11
$ws.Range("E134").Font.Underline = $false
$ws.Range("E134").Font.Color = 0
$ws.Range("E134").NumberFormat = "@"

# --- Row 135: new event (X-BASS @ Purple Pearls, Krefeld) ---
$ws.Range("A135").Value = 45709
$ws.Range("B135").Value = "X-BASS"
$ws.Range("B135").NumberFormat = "@"
$ws.Range("C135").Value = "Purple Pearls"
$ws.Range("C135").NumberFormat = "@"
$ws.Range("D135").Value = "Krefeld"
$ws.Range("D135").NumberFormat = "@"

$ws.Range("E135").Value = "https://www.instagram.com/p/DFnkLImg8RO/?igsh=bXJocHZmMW1zNjho"
$ws.Range("E135").NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Range("E135"), "https://www.instagram.com/p/DFnkLImg8RO/?igsh=bXJocHZmMW1zNjho", "", "", "https://www.instagram.com/p/DFnkLImg8RO/?igsh=bXJocHZmMW1zNjho") | Out-Null

$e135a = $ws.Range("E135").Characters(1, 61)
$e135a.Font.Underline = $true
$e135a.Font.Color = 16711680
$e135b = $ws.Range("E135").Characters(62, 1)
$e135b.Font.Underline = $true
$e135b.Font.Color = 16711680

$ws.Range("E135").Font.Name = "Calibri"
$ws.Range("E135").Font.Size = 11
$ws.Range("E135").Font.Underline = $false
$ws.Range("E135").Font.Color = 0
$ws.Range("E135").NumberFormat = "@"
